# Add custom keyword for jqueryui datepicker
# -> rotate the demo passwords for valid_accounts / invalid_accounts / accounts
#    and leave the "last selection / active sheet" state as it was when the
#    author saved the workbook.

$wb = $excel.ActiveWorkbook

$wsValid   = $wb.Worksheets.Item("valid_accounts")
$wsInvalid = $wb.Worksheets.Item("invalid_accounts")
$wsAll     = $wb.Worksheets.Item("accounts")

# --- new password values -------------------------------------------------
$newPassDemo  = "sPiHQ&YEa6ST``de+"
$newPassTom   = "ok{Ikwnm*wzsaEsD"
$newPassJerry = "gcI#UhR@m(:fsfYU"

# valid_accounts: demo / tom
$wsValid.Range("B2").Value = $newPassDemo
$wsValid.Range("B3").Value = $newPassTom

# invalid_accounts: jerry (bella keeps the stray "jira@2018" value)
$wsInvalid.Range("B2").Value = $newPassJerry

# accounts: demo / tom / jerry
$wsAll.Range("B2").Value = $newPassDemo
$wsAll.Range("B3").Value = $newPassTom
$wsAll.Range("B4").Value = $newPassJerry

# --- restore the per-sheet selections & active tab ------------------------
$wsValid.Range("C16").Select()
$wsAll.Range("F18").Select()

$wsInvalid.Activate()
$wsInvalid.Range("B2").Select()
